$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Text / shared-string content changes -----------------------------
# Order matters: the shared-strings table appends new unique strings in
# the order they are first introduced, and the original workbook appends
# "Documentation", then "To make professional documentation for our
# website.", then the extended w3school/MDN sentence - so we update the
# cells in that same order.
$ws.Range("B24").Value = "Documentation"
$ws.Range("E24").Value = "To make professional documentation for our website."
$ws.Range("D9").Value = "Follow w3school link and HTML/CSS Tutorial on MDN web docs."

# --- Formatting: reuse the existing wrap-only style (no right align) --
# used by F15/F16 for the F17:F20 "Time (hrs)" cells.
[void]$ws.Range("F16").Copy()
[void]$ws.Range("F17:F20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row height for row 10 (skill #4) ----------------------------------
$ws.Range("A10:I10").RowHeight = 94.5

# --- Training hours ("Time (hrs)") column updates ----------------------
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = 2
$ws.Range("F24").Value = 2

# --- View state: zoom to 85% and select F7 as the active cell ----------
[void]$ws.Activate()
[void]$ws.Range("F7").Select()
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
